$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.646.15'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.891.07'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = "'239.68"
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').Value = "'0.4907"
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('D9').Value = "'0.06706"
$ws.Range('D10').Value = '1.892.70'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('D11').Value = "'17.08"
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = "'0.07350"
$ws.Range('E12').Value = '  +1.99%  '
$ws.Range('D13').Value = "'5.153"
$ws.Range('E13').Value = '  +3.51%  '
$ws.Range('D14').Value = "'88.21"
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = "'0.6681"
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '30.589.52'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = "'0.000007876"
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = "'13.41"
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '2.144.64'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = "'5.332"
$ws.Range('E21').Value = '  +12.41%  '
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = "'190.14"
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('E24').Value = '  +3.43%  '
$ws.Range('D25').Value = "'9.533"
$ws.Range('E25').Value = '  +3.35%  '
$ws.Range('D26').Value = "'161.57"
$ws.Range('E26').Value = '  +4.16%  '
$ws.Range('D27').Value = "'18.46"
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = "'1.930"
$ws.Range('E28').Value = '  +4.43%  '
$ws.Range('E29').Value = '  +4.19%  '
$ws.Range('D30').Value = "'4.401"
$ws.Range('E30').Value = '  +4.15%  '
$ws.Range('D31').Value = "'0.09157"
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('E32').Value = '  +3.69%  '
$ws.Range('D33').Value = "'0.05245"
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('D34').Value = "'0.7418"
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('E35').Value = '  +2.33%  '
$ws.Range('D36').Value = "'2.725"
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = "'0.01827"
$ws.Range('E37').Value = '  +0.91%  '
$ws.Range('D38').Value = "'2.689"
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = "'0.9148"
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = "'2.067"
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('D41').Value = "'75.34"
$ws.Range('E41').Value = '  +32.42%  '
$ws.Range('D42').Value = "'0.4421"
$ws.Range('E42').Value = '  +2.54%  '
$ws.Range('D43').Value = "'5.933"
$ws.Range('E43').Value = '  +6.04%  '
$ws.Range('D44').Value = "'106.18"
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('D45').Value = "'0.9930"
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('E46').Value = '  +3.98%  '
$ws.Range('D47').Value = "'7.550"
$ws.Range('E47').Value = '  +3.29%  '
$ws.Range('D48').Value = "'35.49"
$ws.Range('E48').Value = '  +6.73%  '
$ws.Range('D49').Value = "'9.075"
$ws.Range('E49').Value = '  +4.49%  '
$ws.Range('D50').Value = "'0.05839"
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = "'0.3958"
$ws.Range('E51').Value = '  +2.09%  '
